$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I0 / IF columns, rows 2-22
$data = @(
    @(2, 7, 9),
    @(3, 8, 8),
    @(4, 6, 7),
    @(5, 7, 8),
    @(6, 9, 9),
    @(7, 7, 7),
    @(8, 8, 9),
    @(9, 7, 7),
    @(10, 7, 8),
    @(11, 7, 8),
    @(12, 8, 9),
    @(13, 9, 9),
    @(14, 5, 6),
    @(15, 6, 6),
    @(16, 3, 4),
    @(17, 6, 6),
    @(18, 6, 6),
    @(19, 11, 11),
    @(20, 8, 8),
    @(21, 6, 6),
    @(22, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
